$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Formula = '=IF(ISBLANK(B16), "Düsseldorf", B16)'
$ws.Range("A16").Select()
